$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.946.30'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '3.376.28'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.43%  '
$ws.Range("E10").Value = '  -0.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.384'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '3.954.76'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("E13").Value = '  +1.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.70%  '
$ws.Range("D15").Value = '3.379.31'
$ws.Range("E15").Value = '  -0.10%  '
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("D17").Value = '61.049.23'
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("E18").Value = '  -1.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.96%  '
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.547'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").Value = '3.516.65'
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  -2.15%  '
$ws.Range("E27").Value = '  +6.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("E31").Value = '  -0.87%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").Value = '  -4.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.89'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.79'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.32%  '
$ws.Range("D37").Value = '3.412.42'
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.96'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("E39").Value = '  -2.01%  '
$ws.Range("E40").Value = '  -1.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.41%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("E45").Value = '  -2.85%  '
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D47").Value = '2.430.32'
$ws.Range("E47").Value = '  -3.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0260'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.09'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.11%  '
